$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that is bumped by one
# day for every data row (rows 2 through 114) on each automatic update run.
# Old serial 46081 (2026-02-28) -> new serial 46082 (2026-03-01).

for ($row = 2; $row -le 114; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46081) {
        $cell.Value2 = 46082
    }
}
